$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a brand-new paragraph before the "Kerry:" paragraph containing
#    the "Modified: HPActor, dumbledore, HPGridInterface, HPTextInterface
#    and HPWorld classes" line (indented like the other "Created/Modified"
#    paragraphs further down).
# ---------------------------------------------------------------------------
$kerry = $d.Paragraphs.Item(3)
$insertionPoint = $kerry.Range
$insertionPoint.Collapse(1) | Out-Null
$insertionPoint.InsertParagraphBefore() | Out-Null

$newModifiedParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:ind w:firstLine="720"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Modified: </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>HPActor</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve">, </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>dumbl</w:t></w:r>
            <w:r><w:t>ed</w:t></w:r>
            <w:r><w:t>ore</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve">, </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>HP</w:t></w:r>
            <w:r><w:t>GridInterface</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve">, </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>HP</w:t></w:r>
            <w:r><w:t>TextInterface</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> and </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>HPWorld</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> classes</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Paragraphs.Item(3).Range.InsertXML($newModifiedParaXml)

# ---------------------------------------------------------------------------
# 2) The "Kerry: created - immobulus" paragraph (now paragraph 4) keeps its
#    text but gets the document's "_GoBack" bookmark placed right after
#    "immobulus".
# ---------------------------------------------------------------------------
$kerryParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>Kerry:</w:t></w:r>
            <w:r><w:t xml:space="preserve"> created &#8211; </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>immobulus</w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Paragraphs.Item(4).Range.InsertXML($kerryParaXml)

# ---------------------------------------------------------------------------
# 3) The last indented "Modified - GridInterface, TextInterface, HPWorld and
#    Player Classes" paragraph gets its class names prefixed with "HP",
#    extra proofing marks around the comma after "HPTextInterface", and
#    loses the "_GoBack" bookmark (which moved to the Kerry paragraph
#    above).
# ---------------------------------------------------------------------------
$modifiedGridParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:ind w:firstLine="720"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Modified &#8211; </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>HPGridInterface</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve">, </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:t>HPTextInterface</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve">  </w:t></w:r>
            <w:r><w:t>,</w:t></w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>HPWorld</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> and Player Classes </w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$d.Paragraphs.Item(8).Range.InsertXML($modifiedGridParaXml)

Write-Host "Done."
